$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Fix the product name text in both sheets: add missing hyphen after "293"
$ws1.Range("B1").Value = "293-MS-EPP-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"
$ws2.Range("B1").Value = "293-MS-EPP-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"

# Update selection on the input sheet, then switch the active sheet/selection
# to the output sheet (ProductLoanOutput) so it becomes the active tab.
$ws1.Activate()
$ws1.Range("B1").Select()

$ws2.Activate()
$ws2.Range("B1").Select()
